$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(1)

# Shrink the title placeholder's width (flipH="1" xfrm), height stays the same.
$sh.Width = 407.6928346456693

$tr = $sh.TextFrame.TextRange

# Locate "Basics of team development on FPGA" within the full title text
# ("MIPT-MIPS 2013" + line break + the phrase) and replace it with the
# shorter "Basics of teamwork on FPGA", keeping the same start position.
$old = "Basics of team development on FPGA"
$new = "Basics of teamwork on FPGA"
$start = $tr.Text.IndexOf($old) + 1
$whole = $tr.Characters($start, $old.Length)
$whole.Text = $new

# Re-apply the text of each new phrase in place (same total length) so the
# run for each phrase stays distinct instead of being coalesced into one.
$r1 = $tr.Characters($start, 10)
$r1.Text = "Basics of "

$r2 = $tr.Characters($start + 10, 12)
$r2.Text = "teamwork on "

$r3 = $tr.Characters($start + 22, 4)
$r3.Text = "FPGA"
